$d = $word.ActiveDocument
Write-Output "ListTemplates count: $($d.ListTemplates.Count)"
$lt8 = $d.ListTemplates.Item(8)
Write-Output "lt8: $lt8"
$p30 = $d.Paragraphs.Item(30)
$p30.Range.ListFormat.ApplyListTemplateWithLevel($lt8, $false, 1, $false)
Write-Output "applied"
